{"js": "// Append the \"7:\" / \"8:\" work-progress entries to the end of the document,\n// reproducing the exact paragraph/run structure (including the two\n// empty paragraphs, the two en-US \"clamp\" code lines, and the split runs)\n// from the authoritative OOXML diff.\n\n// Raw paragraph OOXML (w:p elements) that must be appended after the\n// document's final paragraph (\"...l\u00e6gge n\u00f8dvendigt slack til.\").\nconst newParagraphsXml =\n  '<w:p/>' +\n  '<w:p><w:r><w:t>7:</w:t></w:r></w:p>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Pr\u00f8vede nu at l\u00e6gge forskellen i gennemsnittene mellem </w:t></w:r>' +\n    '<w:r><w:t>Pressure_tests_Scan_2_40_recon</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> og </w:t></w:r>' +\n    '<w:r><w:t>Pressure_tests_Scan_2_</w:t></w:r>' +\n    '<w:r><w:t>1</w:t></w:r>' +\n    '<w:r><w:t>0_recon</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> til hinanden for at give dem samme v\u00e6rdier. Pr\u00f8vede ogs\u00e5 at clampe p\u00e5 f\u00f8lgende m\u00e5de:</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>image_array[image_array &gt; 40000] = 40000</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>image_array[image_array &lt; 25000] = 25000</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:r><w:t>Igen giver det store negative tal</w:t></w:r>' +\n    '<w:r><w:t>\u2026.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>8</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Pr\u00f8ver at matche p\u00e5 masken nu.</w:t></w:r></w:p>';\n\n// Office.js requires insertOoxml() payloads to be wrapped in the\n// flat-OPC \"pkg:package\" envelope.\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphsXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// Collapsed range at the very end of the document body (after the last\n// paragraph, before the final section break) \u2014 insert our new content there.\nconst endRange = context.document.body.getRange(Word.RangeLocation.end);\nendRange.insertOoxml(flatOpc, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append the \"7:\" / \"8:\" work-progress entries to the end of the document,\n# reproducing the exact paragraph/run structure (including the two empty\n# paragraphs, the two en-US \"clamp\" code lines, and the split runs) from\n# the authoritative OOXML diff.\n\n$d = $word.ActiveDocument\n\n# Raw paragraph OOXML (w:p elements) that must be appended after the\n# document's final paragraph (\"...l\u00e6gge n\u00f8dvendigt slack til.\").\n$newParagraphsXml = '<w:p/>' + `\n  '<w:p><w:r><w:t>7:</w:t></w:r></w:p>' + `\n  '<w:p>' + `\n    '<w:r><w:t xml:space=\"preserve\">Pr\u00f8vede nu at l\u00e6gge forskellen i gennemsnittene mellem </w:t></w:r>' + `\n    '<w:r><w:t>Pressure_tests_Scan_2_40_recon</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> og </w:t></w:r>' + `\n    '<w:r><w:t>Pressure_tests_Scan_2_</w:t></w:r>' + `\n    '<w:r><w:t>1</w:t></w:r>' + `\n    '<w:r><w:t>0_recon</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> til hinanden for at give dem samme v\u00e6rdier. Pr\u00f8vede ogs\u00e5 at clampe p\u00e5 f\u00f8lgende m\u00e5de:</w:t></w:r>' + `\n  '</w:p>' + `\n  '<w:p>' + `\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>image_array[image_array &gt; 40000] = 40000</w:t></w:r>' + `\n  '</w:p>' + `\n  '<w:p>' + `\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>image_array[image_array &lt; 25000] = 25000</w:t></w:r>' + `\n  '</w:p>' + `\n  '<w:p>' + `\n    '<w:r><w:t>Igen giver det store negative tal</w:t></w:r>' + `\n    '<w:r><w:t>\u2026.</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n  '</w:p>' + `\n  '<w:p/>' + `\n  '<w:p/>' + `\n  '<w:p><w:r><w:t>8</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' + `\n  '<w:p><w:r><w:t>Pr\u00f8ver at matche p\u00e5 masken nu.</w:t></w:r></w:p>'\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n        '<w:body>' + $newParagraphsXml + '</w:body>' + `\n      '</w:document>' + `\n    '</pkg:xmlData>' + `\n  '</pkg:part>' + `\n'</pkg:package>'\n\n# Collapse to a zero-length range at the very end of the document body\n# (after the last paragraph, before the final section break) and inject\n# the new paragraphs there.\n$r = $d.Content\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertXML($xml)\n"}
